$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 942.875
$ws.Range("J17").Value = 942.875
$ws.Range("L17").Value = 2828.625
$ws.Range("N17").Value = -3164.625
$ws.Range("H19").Value = 2264.1667
$ws.Range("I19").Value = 2644.3333
$ws.Range("J19").Value = 1123.6666
$ws.Range("K19").Value = 2644.3333
$ws.Range("L19").Value = 1123.6666
$ws.Range("M19").Value = -2469.3333
$ws.Range("N19").Value = -1473.6666
$ws.Range("H34").Value = 8419.857
$ws.Range("I34").Value = 7823.1665
$ws.Range("K34").Value = 7823.1665
$ws.Range("M34").Value = -7620.1665
$ws.Range("H36").Value = 8419.857
$ws.Range("I36").Value = 7823.1665
$ws.Range("K36").Value = 7823.1665
$ws.Range("M36").Value = -7108.1665
$ws.Range("H64").Value = 5497.5
$ws.Range("H67").Value = 5497.5
$ws.Range("H86").Value = 133910.88
$ws.Range("J86").Value = 252376
$ws.Range("L86").Value = 252376
$ws.Range("N86").Value = -254622
$ws.Range("H88").Value = 1474.75
$ws.Range("I88").Value = 1474.75
$ws.Range("K88").Value = 1474.75
$ws.Range("M88").Value = -1068.75
$ws.Range("H89").Value = 133910.88
$ws.Range("J89").Value = 252376
$ws.Range("L89").Value = 1261880
$ws.Range("N89").Value = -1273112
$ws.Range("H91").Value = 1474.75
$ws.Range("I91").Value = 1474.75
$ws.Range("K91").Value = 1474.75
$ws.Range("M91").Value = -70.75
$ws.Range("H98").Value = 920.625
$ws.Range("I98").Value = 809.4286
$ws.Range("K98").Value = 809.4286
$ws.Range("M98").Value = 688.5714
$ws.Range("H116").Value = 7055
$ws.Range("I116").Value = 6666
$ws.Range("K116").Value = 6666
$ws.Range("M116").Value = -3224
$ws.Range("H122").Value = 920.625
$ws.Range("I122").Value = 809.4286
$ws.Range("K122").Value = 2428.2858
$ws.Range("M122").Value = 21.71420000000035

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1946.1818
$ws.Range("I20").Value = 1801.7142
$ws.Range("J20").Value = 2199
$ws.Range("K20").Value = 1801.7142
$ws.Range("L20").Value = 2199
$ws.Range("M20").Value = -1554.7142
$ws.Range("N20").Value = -2693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12973.5
$ws.Range("I86").Value = 13290
$ws.Range("K86").Value = 13290
$ws.Range("M86").Value = -12167
$ws.Range("H89").Value = 12973.5
$ws.Range("I89").Value = 13290
$ws.Range("K89").Value = 66450
$ws.Range("M89").Value = -60834
$ws.Range("H134").Value = 2276.724
$ws.Range("I134").Value = 2258.68
$ws.Range("J134").Value = 2389.5
$ws.Range("K134").Value = 6776.039999999999
$ws.Range("L134").Value = 7168.5
$ws.Range("M134").Value = -4241.039999999999
$ws.Range("N134").Value = -12238.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 208.36363
$ws.Range("I6").Value = 219.2
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 657.5999999999999
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -544.5999999999999
$ws.Range("N6").Value = -526
$ws.Range("H103").Value = 818.2222
$ws.Range("I103").Value = 40.666668
$ws.Range("J103").Value = 1207
$ws.Range("K103").Value = 122.000004
$ws.Range("L103").Value = 3621
$ws.Range("M103").Value = 756.999996
$ws.Range("N103").Value = -5379
$ws.Range("H106").Value = 19124.875
$ws.Range("J106").Value = 19799.8
$ws.Range("L106").Value = 59399.39999999999
$ws.Range("N106").Value = -61291.39999999999
$ws.Range("H128").Value = 580237
$ws.Range("I128").Value = 580237
$ws.Range("K128").Value = 1740711
$ws.Range("M128").Value = -1735731
$ws.Range("H131").Value = 590722.1
$ws.Range("J131").Value = 627579.75
$ws.Range("L131").Value = 1882739.25
$ws.Range("N131").Value = -1892819.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 2962.6667
$ws.Range("J80").Value = 2962.6667
$ws.Range("L80").Value = 2962.6667
$ws.Range("N80").Value = -4958.6667
$ws.Range("H83").Value = 2962.6667
$ws.Range("J83").Value = 2962.6667
$ws.Range("L83").Value = 14813.3335
$ws.Range("N83").Value = -24797.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1829.2727
$ws.Range("I22").Value = 1494.1428
$ws.Range("K22").Value = 1494.1428
$ws.Range("M22").Value = -1199.1428
$ws.Range("H27").Value = 1829.2727
$ws.Range("I27").Value = 1494.1428
$ws.Range("K27").Value = 1494.1428
$ws.Range("M27").Value = -1387.1428
$ws.Range("H63").Value = 83999
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 83999
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 80732
$ws.Range("I74").Value = 81098.5
$ws.Range("K74").Value = 81098.5
$ws.Range("M74").Value = -80100.5
$ws.Range("H77").Value = 80732
$ws.Range("I77").Value = 81098.5
$ws.Range("K77").Value = 243295.5
$ws.Range("M77").Value = -238303.5
$ws.Range("H120").Value = 43000
$ws.Range("J120").Value = 43000
$ws.Range("L120").Value = 43000
$ws.Range("N120").Value = -52676
$ws.Range("H132").Value = 8327.200000000001
$ws.Range("I132").Value = 10310.429
$ws.Range("K132").Value = 30931.287
$ws.Range("M132").Value = -28401.287
$ws.Range("H136").Value = 5196.2
$ws.Range("I136").Value = 4499.5
$ws.Range("K136").Value = 13498.5
$ws.Range("M136").Value = -10948.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 99999
$ws.Range("J2").Value = 99999
$ws.Range("L2").Value = 99999
$ws.Range("N2").Value = -100223
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4888
$ws.Range("H75").Value = 73203.25
$ws.Range("I75").Value = 69905
$ws.Range("J75").Value = 74302.664
$ws.Range("K75").Value = 69905
$ws.Range("L75").Value = 74302.664
$ws.Range("M75").Value = -68969
$ws.Range("N75").Value = -76174.664
$ws.Range("H78").Value = 73203.25
$ws.Range("I78").Value = 69905
$ws.Range("J78").Value = 74302.664
$ws.Range("K78").Value = 209715
$ws.Range("L78").Value = 222907.992
$ws.Range("M78").Value = -205035
$ws.Range("N78").Value = -232267.992
$ws.Range("H136").Value = 3579.1538
$ws.Range("I136").Value = 4183.1
$ws.Range("K136").Value = 12549.3
$ws.Range("M136").Value = -9999.300000000001
